# 2017-7-4: add 3 new coll and their related rule and test case
#
# The sheet already has the "currentColl / currentPage / recIdArr / searchParams /
# filterFieldValue / editSubField / eventField" header row; the only real content
# change is a typo fix in the very first header cell (B1): "recorderInfo" -> "recordInfo".
# The view also scrolls so that column F is the left-most visible column and the
# active selection becomes E1 (instead of the scrolled-down/right A7 / E17 view).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the typo in the header cell B1.
$ws.Range("B1").Value = "recordInfo"

# Update the view: scroll so column F is left-most and select E1.
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1
$ws.Range("E1").Select()
